$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows: one at row 2, one at row 5 (post first insert) to make room
# for the two brand-new Diot entries described by the diff.
$ws.Rows.Item(2).Insert() | Out-Null
$ws.Rows.Item(5).Insert() | Out-Null

# Folio (column G) values are numeric-looking text with a trailing space in the
# source data; force text format so COM does not coerce them into real numbers.
$ws.Range("G2:G15").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Diot"
$ws.Range("B2").Value = "AGR1606133D8 "
$ws.Range("C2").Value = "AGR1606133D80DOTAAN1NCN5J14451.dec "
$ws.Range("D2").Value = "Tamaño:1530 "
$ws.Range("E2").Value = "19/05/2023 "
$ws.Range("F2").Value = "09:42:49 "
$ws.Range("G2").Value = "387312635 "
$ws.Range("H2").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\DIOT AGROL ABRIL 2023   L.pdf"

# Row 3
$ws.Range("A3").Value = "Diot"
$ws.Range("B3").Value = "AGR1606133D8 "
$ws.Range("C3").Value = "AGR1606133D80DOTAAN1NCN2S18121.dec "
$ws.Range("D3").Value = "Tamaño:1474 "
$ws.Range("E3").Value = "06/03/2023 "
$ws.Range("F3").Value = "18:35:18 "
$ws.Range("G3").Value = "383452104 "
$ws.Range("H3").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot agrol MARZO 2023.pdf"

# Row 4
$ws.Range("A4").Value = "Diot"
$ws.Range("B4").Value = "BOD1702215A0 "
$ws.Range("C4").Value = "BOD1702215A00DOTAAN1NCN2R22101.dec "
$ws.Range("D4").Value = "Tamaño:3442 "
$ws.Range("E4").Value = "28/02/2023 "
$ws.Range("F4").Value = "23:04:18 "
$ws.Range("G4").Value = "383145045 "
$ws.Range("H4").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot bricks febrero 2023.pdf"

# Row 5
$ws.Range("A5").Value = "Diot"
$ws.Range("B5").Value = "BOD1702215A0 "
$ws.Range("C5").Value = "BOD1702215A00DOTAAN1NCN5T16081.dec "
$ws.Range("D5").Value = "Tamaño:2946 "
$ws.Range("E5").Value = "30/05/2023 "
$ws.Range("F5").Value = "19:14:29 "
$ws.Range("G5").Value = "388083343 "
$ws.Range("H5").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\DIOT BRICKS OBRAS Y DESARROLLO ABRIL 2023.pdf"

# Row 6
$ws.Range("A6").Value = "Diot"
$ws.Range("B6").Value = "CAP990628715 "
$ws.Range("C6").Value = "CAP9906287150DOTAAN1NCN2401291.dec "
$ws.Range("D6").Value = "Tamaño:1482 "
$ws.Range("E6").Value = "15/02/2023 "
$ws.Range("F6").Value = "16:45:04 "
$ws.Range("G6").Value = "382257225 "
$ws.Range("H6").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot cuadras asesores DICIEMBRE 2022.pdf"

# Row 7
$ws.Range("A7").Value = "Diot"
$ws.Range("B7").Value = "DES160421EU7 "
$ws.Range("C7").Value = "DES160421EU70DOTAAN1NCN2O15421.dec "
$ws.Range("D7").Value = "Tamaño:2834 "
$ws.Range("E7").Value = "27/02/2023 "
$ws.Range("F7").Value = "09:55:45 "
$ws.Range("G7").Value = "382924957 "
$ws.Range("H7").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot desoflex 27-02-23.pdf"

# Row 8
$ws.Range("A8").Value = "Diot"
$ws.Range("B8").Value = "DLS1403063Z9 "
$ws.Range("C8").Value = "DLS1403063Z90DOTAAN1NCN2O22481.dec "
$ws.Range("D8").Value = "Tamaño:1378 "
$ws.Range("E8").Value = "28/02/2023 "
$ws.Range("F8").Value = "16:42:07 "
$ws.Range("G8").Value = "383101606 "
$ws.Range("H8").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot diseños luna 28-02-23.pdf"

# Row 9
$ws.Range("A9").Value = "Diot"
$ws.Range("B9").Value = "JARH570121P36 "
$ws.Range("D9").Value = "Tamaño:1130 "
$ws.Range("E9").Value = "28/02/2023 "
$ws.Range("F9").Value = "16:28:11 "
$ws.Range("G9").Value = "383098681 "
$ws.Range("H9").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot hector jauregui rios  28-02-23.pdf"

# Row 10
$ws.Range("A10").Value = "Diot"
$ws.Range("H10").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot hector jauregui rios 14-02-23.pdf"

# Row 11
$ws.Range("A11").Value = "Diot"
$ws.Range("B11").Value = "INT0901197U0 "
$ws.Range("C11").Value = "INT0901197U00DOTAAN1NCN3H23261.dec "
$ws.Range("D11").Value = "Tamaño:1354 "
$ws.Range("E11").Value = "21/03/2023 "
$ws.Range("F11").Value = "09:11:45 "
$ws.Range("G11").Value = "384026542 "
$ws.Range("H11").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\DIOT INTEREXPORTA ENERO 2023.pdf"

# Row 12
$ws.Range("A12").Value = "Diot"
$ws.Range("B12").Value = "INT0901197U0 "
$ws.Range("C12").Value = "INT0901197U00DOTAAN1NCN4422551.dec "
$ws.Range("D12").Value = "Tamaño:1170 "
$ws.Range("E12").Value = "25/04/2023 "
$ws.Range("F12").Value = "12:15:00 "
$ws.Range("G12").Value = "385804991 "
$ws.Range("H12").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\DIOT INTEREXPORTA FEBRERO 2023.pdf"

# Row 13
$ws.Range("A13").Value = "Diot"
$ws.Range("B13").Value = "ZALM740419353 "
$ws.Range("D13").Value = "Tamaño:2250 "
$ws.Range("E13").Value = "13/03/2023 "
$ws.Range("F13").Value = "10:38:39 "
$ws.Range("G13").Value = "383699154 "
$ws.Range("H13").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot monica alejandra zarate losa 13-03-23.pdf"

# Row 14
$ws.Range("A14").Value = "Diot"
$ws.Range("B14").Value = "IAMP661003NP2 "
$ws.Range("D14").Value = "Tamaño:1218 "
$ws.Range("E14").Value = "28/02/2023 "
$ws.Range("F14").Value = "16:37:13 "
$ws.Range("G14").Value = "383100589 "
$ws.Range("H14").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot pascual ibarra 28-02-23.pdf"

# Row 15
$ws.Range("A15").Value = "Diot"
$ws.Range("B15").Value = "MOAR741018D36 "
$ws.Range("D15").Value = "Tamaño:5946 "
$ws.Range("E15").Value = "01/03/2023 "
$ws.Range("F15").Value = "08:34:02 "
$ws.Range("G15").Value = "383152149 "
$ws.Range("H15").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\diot rigoberto mora 01-03-23.pdf"
